$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 (VT200-0851 / "Signal showIcon all property"): add icon-position
# validations to validate4 in column H, and grow the row to fit the extra
# lines.
$ws.Cells.Item(3, 8).Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signal JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0851`n};`nvalidate4`n{`nvalidate_Screenshot=VT200-0851`nvalidate_Iconposition=signalview_xpath,left,20`nvalidate_Iconposition=signalview_xpath,top,40`n};"
$ws.Rows.Item(3).RowHeight = 203.25

# Row 6 (VT200-0854 / "Signal showIcon left with 80"): drop the screenshot
# step from the script and replace the screenshot validation with an icon
# position check.
$ws.Cells.Item(6, 7).Value = "wait(3);`nvalidate1;`nlink_Click(signal_test_link);`nvalidate2;`nSelectTestToRun(VT200_0854_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nvalidate4;"
$ws.Cells.Item(6, 8).Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signal JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0854`n};`nvalidate4`n{`nvalidate_Iconposition=signalview_xpath,left,40`n};"

# Row 7 (VT200-0855 / "Signal showIcon top with 80"): same shape of change
# as row 6, but for the "top" icon position.
$ws.Cells.Item(7, 7).Value = "wait(3);`nvalidate1;`nlink_Click(signal_test_link);`nvalidate2;`nSelectTestToRun(VT200_0855_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nvalidate4;"
$ws.Cells.Item(7, 8).Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signal JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0855`n};`nvalidate4`n{`nvalidate_Iconposition=signalview_xpath,top,40`n};"

# Row 13 (VT200-0861 / "Signal hideIcon after calling showIcon method"):
# drop the second screenshot step, and swap both screenshot validations for
# isIconDisplayed checks.
$ws.Cells.Item(13, 7).Value = "wait(3);`nvalidate1;`nlink_Click(signal_test_link);`nvalidate2;`nSelectTestToRun(VT200_0861_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nTakeScreenshot(VT200-0861-01);`nwait(2);`nvalidate4;`nwait(12);`nvalidate5;"
$ws.Cells.Item(13, 8).Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signal JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0861`n};`nvalidate4`n{`nvalidate_isIconDisplayed=signalview_xpath,true`n};`nvalidate5`n{`nvalidate_isIconDisplayed=signalview_xpath,false`n};"

# Select G1 on the active sheet, matching the saved selection state.
$ws.Range("G1").Select()
